# Append two new test-case rows to the ScriptMasterSheet (TestScriptMasterSheet):
#   testT4763  | 1 | YES
#   testT4275  | 1 | YES
# These follow the same three-column layout (TESTCASE, SCRIPT_ITERATION,
# EXECUTE_FLAG) as every other row in the table, directly below the last
# existing row (row 32, NonInvestigativeCaseDataSetup1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the last populated row in column A (mirrors the user pressing
# Ctrl+End / typing below the last existing table row).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

$newRow1 = $lastRow + 1
$newRow2 = $lastRow + 2

$ws.Cells.Item($newRow1, 1).Value = "testT4763"
$ws.Cells.Item($newRow1, 2).Value = 1
$ws.Cells.Item($newRow1, 3).Value = "YES"

$ws.Cells.Item($newRow2, 1).Value = "testT4275"
$ws.Cells.Item($newRow2, 2).Value = 1
$ws.Cells.Item($newRow2, 3).Value = "YES"

# Update the active selection to mirror the author's final view of the
# sheet after adding the new rows.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 18
$ws.Range("F29").Select()
